$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Restricciones_del_follower
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")

# Pre-format the numeric-looking text cells as Text ("@") so that assigning a
# numeric-looking string (e.g. "-2.49") keeps it stored as a text value
# (shared string) instead of being auto-converted into a real number by
# Excel's input parser.
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("D2:F2").NumberFormat = "@"
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("D3:F3").NumberFormat = "@"
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("D4:F4").NumberFormat = "@"
$ws3.Range("B5").NumberFormat = "@"
$ws3.Range("D5:F5").NumberFormat = "@"
$ws3.Range("B6").NumberFormat = "@"
$ws3.Range("D6:F6").NumberFormat = "@"

$ws3.Range("A2").Value = "4.49 - x - 0.5y"
$ws3.Range("B2").Value = "-2.49"
$ws3.Range("C2").Value = "J_0_L0_v"
$ws3.Range("D2").Value = "0.62"
$ws3.Range("E2").Value = "7.6"
$ws3.Range("F2").Value = "0"

$ws3.Range("A3").Value = "-4.4125 - 0.25x + y"
$ws3.Range("B3").Value = "2.4124999999999996"
$ws3.Range("C3").Value = "J_0_L0_v"
$ws3.Range("D3").Value = "0.96"
$ws3.Range("E3").Value = "9.399999999999999"
$ws3.Range("F3").Value = "1.6"

$ws3.Range("A4").Value = "-4.49 + x + 0.5y"
$ws3.Range("B4").Value = "-3.51"
$ws3.Range("C4").Value = "J_0_LP_v"
$ws3.Range("D4").Value = "0.88"
$ws3.Range("E4").Value = "-2.3000000000000003"
$ws3.Range("F4").Value = "-3.7"

$ws3.Range("A5").Value = "-11.809999999999999 + x - 2y"
$ws3.Range("B5").Value = "-9.809999999999999"
$ws3.Range("C5").Value = "J_Ne_L0_v"
$ws3.Range("D5").Value = "0.29"
$ws3.Range("E5").Value = "5.4"
$ws3.Range("F5").Value = "0"

$ws3.Range("A6").Value = "-4.92 - y"
$ws3.Range("B6").Value = "-4.92"
$ws3.Range("C6").Value = "J_Ne_L0_v"
$ws3.Range("D6").Value = "0.34"
$ws3.Range("E6").Value = "0"
$ws3.Range("F6").Value = "0.2"

# ---------------------------------------------------------------------------
# Sheet: Punto_modificado
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Punto_modificado")
$ws4.Range("A2:B2").NumberFormat = "@"
$ws4.Range("A2").Value = "2.0300000000000002"
$ws4.Range("B2").Value = "4.92"

# ---------------------------------------------------------------------------
# Sheet: Vector_bf
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Vector_bf")
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "0.8300000000000001"

# ---------------------------------------------------------------------------
# Sheet: Vector_BF
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Vector_BF")
$ws6.Range("A2:A3").NumberFormat = "@"
$ws6.Range("A2").Value = "5.85"
$ws6.Range("A3").Value = "5.350000000000002"
